$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) column updates
# Cells whose new values are plain decimal numbers must be forced to Text format
# first, so Excel does not auto-convert the assigned string into a numeric value
# (the source data is stored as text, matching the original inline-string cells).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.41"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4541"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3525"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07364"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.074"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.37"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.906"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.052"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.88"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001052"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.57"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.746"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.071"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.34"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.91"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.042"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.79"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.043"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09104"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.657"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02262"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05946"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2048"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6219"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.863"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.372"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.683"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.99"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.695"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5775"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.65"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.921"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.110"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.93"

# Cells whose new values contain multiple dots (e.g. "27.437.53") are never
# auto-converted to numbers by Excel, so no text-format override is required.
$ws.Range("D2").Value = "27.437.53"
$ws.Range("D3").Value = "1.735.84"
$ws.Range("D16").Value = "1.728.96"
$ws.Range("D23").Value = "27.467.52"
$ws.Range("D28").Value = "1.926.92"

# Volume(1h) (E) column updates
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +6.98%  "
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("E16").Value = "  -3.14%  "
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("E29").Value = "  -3.85%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  -4.86%  "
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("E50").Value = "  -5.01%  "
$ws.Range("E51").Value = "  -3.34%  "
